$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246 - this shifts the former rows 246..284
# down to 247..285, matching the rest of the data set already in place.
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row 246 with its data (a new weekly price
# observation for Pina - Caramelo - Segunda at Vega Modelo de Temuco).
$ws.Range("A246").Value = 10
$ws.Range("B246").Value = "Vega Modelo de Temuco"
$ws.Range("C246").Value = "La Araucanía"
$ws.Range("D246").Value = 44474
$ws.Range("D246").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E246").Value = 9
$ws.Range("F246").Value = "Fruta"
$ws.Range("G246").Value = 100108
$ws.Range("H246").Value = "Tropicales y subtropicales"
$ws.Range("I246").Value = 100108005
$ws.Range("J246").Value = "Piña"
$ws.Range("K246").Value = "Caramelo"
$ws.Range("L246").Value = "Segunda"
$ws.Range("M246").Value = 30
$ws.Range("N246").Value = 20000
$ws.Range("O246").Value = 20000
$ws.Range("P246").Value = 20000
$ws.Range("Q246").Value = "$/caja 14 unidades"
$ws.Range("R246").Value = "Ecuador"
$ws.Range("S246").Value = 1429
$ws.Range("T246").Value = 14
